# Re-word four of the Likert-style prompt strings in column B of Sheet1.
# (Same row/position, just updated wording -- Excel will drop the old
# shared-string entries that become unused and append the new text at the
# end of the shared-strings table.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value  = "I have almost never experienced such issues in the homes where I've stayed."
$ws.Range("B7").Value  = "In the homes where I've lived, I have experienced all the issues mentioned in the story."
$ws.Range("B12").Value = "There have been no issues."
$ws.Range("B16").Value = "I have almost never experienced temperature problems in the homes where I have been."

# Move the selection/view to match the saved workbook state (B16 selected,
# scrolled so column B is the left-most visible column).
$ws.Range("B16").Select()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
